$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'59.273.12"
$ws.Range("E2").Value = "  +3.41%  "

# Row 3
$ws.Range("D3").Value = "'2.544.18"
$ws.Range("E3").Value = "  +5.55%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "'526.47"
$ws.Range("E5").Value = "  +3.19%  "

# Row 6
$ws.Range("D6").Value = "'134.94"
$ws.Range("E6").Value = "  +4.80%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.18%  "

# Row 8
$ws.Range("D8").Value = "'0.568"
$ws.Range("E8").Value = "  +4.01%  "

# Row 9
$ws.Range("D9").Value = "'2.542.70"
$ws.Range("E9").Value = "  +5.32%  "

# Row 10
$ws.Range("D10").Value = "'0.0991"
$ws.Range("E10").Value = "  +4.09%  "

# Row 11
$ws.Range("E11").Value = "  -0.80%  "

# Row 12
$ws.Range("E12").Value = "  +0.97%  "

# Row 13
$ws.Range("E13").Value = "  +1.75%  "

# Row 14
$ws.Range("D14").Value = "'2.995.81"
$ws.Range("E14").Value = "  +5.56%  "

# Row 15
$ws.Range("D15").Value = "'59.206.20"
$ws.Range("E15").Value = "  +3.39%  "

# Row 16
$ws.Range("D16").Value = "'22.43"
$ws.Range("E16").Value = "  +5.19%  "

# Row 17
$ws.Range("E17").Value = "  +3.83%  "

# Row 18
$ws.Range("D18").Value = "'2.543.98"
$ws.Range("E18").Value = "  +5.28%  "

# Row 19
$ws.Range("D19").Value = "'10.75"
$ws.Range("E19").Value = "  +4.28%  "

# Row 20
$ws.Range("D20").Value = "'323.83"
$ws.Range("E20").Value = "  +3.41%  "

# Row 21
$ws.Range("E21").Value = "  +3.58%  "

# Row 22
$ws.Range("D22").Value = "'6.16"
$ws.Range("E22").Value = "  +10.07%  "

# Row 23
$ws.Range("E23").Value = "  +0.21%  "

# Row 24
$ws.Range("D24").Value = "'65.42"
$ws.Range("E24").Value = "  +3.24%  "

# Row 25
$ws.Range("D25").Value = "'0.412"
$ws.Range("E25").Value = "  +3.02%  "

# Row 26
$ws.Range("E26").Value = "  -0.03%  "

# Row 27
$ws.Range("E27").Value = "  +1.69%  "

# Row 28
$ws.Range("D28").Value = "'7.53"
$ws.Range("E28").Value = "  +5.38%  "

# Row 29
$ws.Range("E29").Value = "  +6.22%  "

# Row 30
$ws.Range("E30").Value = "  +7.89%  "

# Row 31
$ws.Range("E31").Value = "  +5.25%  "

# Row 32
$ws.Range("D32").Value = "'169.64"
$ws.Range("E32").Value = "  +0.49%  "

# Row 33
$ws.Range("E33").Value = "  +3.72%  "

# Row 34
$ws.Range("E34").Value = "  -0.03%  "

# Row 35
$ws.Range("E35").Value = "  +0.02%  "

# Row 36
$ws.Range("E36").Value = "  +3.70%  "

# Row 37
$ws.Range("E37").Value = "  +0.64%  "

# Row 38
$ws.Range("D38").Value = "'3.99"
$ws.Range("E38").Value = "  +4.10%  "

# Row 39
$ws.Range("D39").Value = "'1.52"
$ws.Range("E39").Value = "  +6.26%  "

# Row 40
$ws.Range("D40").Value = "'36.79"
$ws.Range("E40").Value = "  +1.47%  "

# Row 41
$ws.Range("D41").Value = "'0.788"
$ws.Range("E41").Value = "  +2.80%  "

# Row 42
$ws.Range("D42").Value = "'280.11"
$ws.Range("E42").Value = "  +6.08%  "

# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.49"
$ws.Range("E43").Value = "  +4.42%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'133.84"
$ws.Range("E44").Value = "  +10.41%  "

# Row 45
$ws.Range("D45").Value = "'5.10"
$ws.Range("E45").Value = "  +4.97%  "

# Row 46
$ws.Range("E46").Value = "  +4.32%  "

# Row 47
$ws.Range("D47").Value = "'0.0924"
$ws.Range("E47").Value = "  +3.02%  "

# Row 48
$ws.Range("E48").Value = "  +6.25%  "

# Row 49
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0218"
$ws.Range("E49").Value = "  +4.41%  "

# Row 50
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'17.20"
$ws.Range("E50").Value = "  +4.93%  "

# Row 51
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "'1.760.72"
$ws.Range("E51").Value = "  +4.57%  "
